# Update cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.114.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.48%  '

$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5099'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3903'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09747'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +21.46%  '

$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("E11").Value = '  -1.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.433'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9992'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.815.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.334'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001139'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("E22").Value = '  +1.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.117.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.68%  '

$ws.Range("E24").Value = '  -0.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.220'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.426'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.014.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1089'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '

$ws.Range("E32").Value = '  -0.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.624'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.09%  '

$ws.Range("E34").Value = '  -1.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06875'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.73%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.061'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02324'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2167'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.29%  '

$ws.Range("E39").Value = '  -7.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.017'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6168'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9995'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.151'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.51%  '

$ws.Range("E44").Value = '  +0.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5939'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.287'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.699'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.958'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.182'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06767'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.76%  '
